# ---------------------------------------------------------------------------
# Applies the "Added calculations and evaluation for corpus 3" commit:
#   * renames Sheet1 -> Statistics
#   * adds two new sheets: Evalua, Legend
#   * inserts a new (blank) column G in Statistics, shifting old G:I to H:J
#   * appends 6 new data rows (8-13) to Statistics for corpus "3"
#   * populates Evalua with the evaluation-metric header row (+ hyperlink)
#   * populates Legend with the glossary / legend table (+ hyperlink)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheet and add the two new ones in order.
# ---------------------------------------------------------------------------
$wsStats = $wb.Worksheets.Item(1)
$wsStats.Name = "Statistics"

$wsEval = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsStats)
$wsEval.Name = "Evalua"

$wsLegend = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsEval)
$wsLegend.Name = "Legend"

# ---------------------------------------------------------------------------
# 2. Statistics sheet: insert a new blank column at G (old G/H/I -> H/I/J)
# ---------------------------------------------------------------------------
$wsStats.Columns("G").Insert()

# ---------------------------------------------------------------------------
# 3. Statistics sheet: append the corpus-3 data block (rows 8-13)
# ---------------------------------------------------------------------------

# Row 8 - avg
$wsStats.Cells.Item(8, 1).Value2 = "3(with negatives)"
$wsStats.Cells.Item(8, 2).Value2 = "avg"
$wsStats.Cells.Item(8, 3).Value2 = 0.8499
$wsStats.Cells.Item(8, 4).Value2 = 0.5292
$wsStats.Cells.Item(8, 5).Value2 = 0.1023
$wsStats.Cells.Item(8, 6).Value2 = 0.2914
$wsStats.Cells.Item(8, 7).Value2 = 0.6194
$wsStats.Cells.Item(8, 8).Value2 = 0.4707
$wsStats.Cells.Item(8, 9).Value2 = 0.0912
$wsStats.Cells.Item(8, 10).Value2 = 0.8965

# Row 9 - min
$wsStats.Cells.Item(9, 2).Value2 = "min"
$wsStats.Cells.Item(9, 3).Value2 = 0
$wsStats.Cells.Item(9, 4).Value2 = 0
$wsStats.Cells.Item(9, 5).Value2 = -0.17
$wsStats.Cells.Item(9, 6).Value2 = -0.1088
$wsStats.Cells.Item(9, 7).Value2 = 0
$wsStats.Cells.Item(9, 8).Value2 = 0.1006
$wsStats.Cells.Item(9, 9).Value2 = -0.17
$wsStats.Cells.Item(9, 10).Value2 = 0.3725

# Row 10 - max
$wsStats.Cells.Item(10, 2).Value2 = "max"
$wsStats.Cells.Item(10, 3).Value2 = 1
$wsStats.Cells.Item(10, 4).Value2 = 1
$wsStats.Cells.Item(10, 5).Value2 = 0.7667
$wsStats.Cells.Item(10, 6).Value2 = 0.9209
$wsStats.Cells.Item(10, 7).Value2 = 1
$wsStats.Cells.Item(10, 8).Value2 = 1
$wsStats.Cells.Item(10, 9).Value2 = 1
$wsStats.Cells.Item(10, 10).Value2 = 1

# Row 11 - avg
$wsStats.Cells.Item(11, 1).Value2 = "3 (only positives)"
$wsStats.Cells.Item(11, 2).Value2 = "avg"
$wsStats.Cells.Item(11, 3).Value2 = 0.9804
$wsStats.Cells.Item(11, 4).Value2 = 0.5347
$wsStats.Cells.Item(11, 5).Value2 = 0.1566
$wsStats.Cells.Item(11, 6).Value2 = 0.415
$wsStats.Cells.Item(11, 7).Value2 = 0.9491
$wsStats.Cells.Item(11, 8).Value2 = 0.572
$wsStats.Cells.Item(11, 9).Value2 = 0.144
$wsStats.Cells.Item(11, 10).Value2 = 0.9928

# Row 12 - min
$wsStats.Cells.Item(12, 2).Value2 = "min"
$wsStats.Cells.Item(12, 3).Value2 = 0.381
$wsStats.Cells.Item(12, 4).Value2 = 0
$wsStats.Cells.Item(12, 5).Value2 = -0.1245
$wsStats.Cells.Item(12, 6).Value2 = -0.0196
$wsStats.Cells.Item(12, 7).Value2 = 0
$wsStats.Cells.Item(12, 8).Value2 = 0.3121
$wsStats.Cells.Item(12, 9).Value2 = -0.1245
$wsStats.Cells.Item(12, 10).Value2 = 0.4991

# Row 13 - max
$wsStats.Cells.Item(13, 2).Value2 = "max"
$wsStats.Cells.Item(13, 3).Value2 = 1
$wsStats.Cells.Item(13, 4).Value2 = 1
$wsStats.Cells.Item(13, 5).Value2 = 0.7667
$wsStats.Cells.Item(13, 6).Value2 = 0.9209
$wsStats.Cells.Item(13, 7).Value2 = 1
$wsStats.Cells.Item(13, 8).Value2 = 0.8371
$wsStats.Cells.Item(13, 9).Value2 = 0.6511
$wsStats.Cells.Item(13, 10).Value2 = 1

# Keep the same active-cell / selection style the author ended up with.
$wsStats.Range("F16").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Evalua sheet: evaluation-metric header row
# ---------------------------------------------------------------------------
$wsEval.Range("A1").Value2 = "Version"
$wsEval.Range("B1").Value2 = "MRR@1"
$wsEval.Range("C1").Value2 = "R-Precision"
$wsEval.Range("D1").Value2 = "Full MRR"

$wsEval.Hyperlinks.Add($wsEval.Range("B1"), "https://en.wikipedia.org/wiki/Evaluation_measures_(information_retrieval)", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "MRR@1") | Out-Null

$wsEval.Columns("C").AutoFit()

$wsEval.Range("B1:D1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Legend sheet: glossary table describing the Statistics/Evalua columns
# ---------------------------------------------------------------------------
$wsLegend.Range("A1").Value2 = "Statistics "

$wsLegend.Range("B2").Value2 = "Version"
$wsLegend.Range("C2").Value2 = "Version of the corpus"

$wsLegend.Range("B3").Value2 = "name"
$wsLegend.Range("C3").Value2 = "calculated similarities of names"

$wsLegend.Range("B4").Value2 = "author"
$wsLegend.Range("C4").Value2 = "calculated similarities of authors"

$wsLegend.Range("B5").Value2 = "keywords"
$wsLegend.Range("C5").Value2 = "calculated similarities of keywordss"

$wsLegend.Range("B6").Value2 = "paragraph"
$wsLegend.Range("C6").Value2 = "calculated similarities of paragraphs"

$wsLegend.Range("B7").Value2 = "average"
$wsLegend.Range("C7").Value2 = "averaged similarities per row"

$wsLegend.Range("B8").Value2 = "min"
$wsLegend.Range("C8").Value2 = "maximum similarity per row"

$wsLegend.Range("B9").Value2 = "max"
$wsLegend.Range("C9").Value2 = "minimum similarity per row"

$wsLegend.Range("A10").Value2 = "Evaluation"
$wsLegend.Range("B10").Value2 = "MRR@1"
$wsLegend.Range("C10").Value2 = "the proportion of cases where the correct item appears in the very first position"

$wsLegend.Range("B11").Value2 = "R-Precision"
$wsLegend.Range("C11").Value2 = "the fraction of all true items that are retrieved within the top R results, where R is the number of true items for that query"

$wsLegend.Range("B12").Value2 = "Full MRR"
$wsLegend.Range("C12").Value2 = "the average over all queries of 1 divided by the rank of the first correct item"

$wsLegend.Hyperlinks.Add($wsLegend.Range("B10"), "https://en.wikipedia.org/wiki/Evaluation_measures_(information_retrieval)", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "MRR@1") | Out-Null

$wsLegend.Range("C15").Select() | Out-Null

# ---------------------------------------------------------------------------
# Leave the Statistics tab active/selected, matching the source workbook.
# ---------------------------------------------------------------------------
$wsStats.Activate() | Out-Null
